$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the column header in B1 from "Population" to "population_name"
$ws.Range("B1").Value = "population_name"

# Move the active cell selection to B1 (matches the saved selection state in the diff)
$ws.Range("B1").Select()
